$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetings")

# Update header dates for I6 and J6 (previously "TBC" placeholders)
$ws.Range("I6").NumberFormat = "d-mmm"
$ws.Range("I6").Value = 43172
$ws.Range("J6").NumberFormat = "d-mmm"
$ws.Range("J6").Value = 43174

# Fill in attendance values for the new meeting columns (I and J) for rows 7-16
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1

$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1

$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1

$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1

$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0

$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1

$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1

# Move the active selection to J7 as in the final workbook
$ws.Range("J7").Select()
